# Automatic update of files.
#
# The rows 2-16 of the active sheet each describe one observation record
# (Id, Taxonsorteringsordning, Valideringsstatus, Rödlistade, TaxonId,
# Artnamn, Vetenskapligt namn, Auktor, ..., Ost, Nord, ...).
# The edit re-associates the per-record data held in columns
# A, B, D, E, F, G, H, Q, R between rows, leaving every other column
# (which was already identical across the affected rows) untouched.
#
# Concretely, the new content of each row below is taken from the row
# indicated in the map (read as "row -> row its A/B/D/E/F/G/H/Q/R values
# come from"):
#   2 <- 4, 3 <- 5, 4 <- 2, 5 <- 6, 6 <- 7, 7 <- 8, 8 <- 9, 9 <- 10,
#   10 <- 11, 11 <- 12, 12 <- 13, 13 <- 14, 14 <- 15, 15 <- 16, 16 <- 3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sourceRows = @{
    2  = 4
    3  = 5
    4  = 2
    5  = 6
    6  = 7
    7  = 8
    8  = 9
    9  = 10
    10 = 11
    11 = 12
    12 = 13
    13 = 14
    14 = 15
    15 = 16
    16 = 3
}

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

# Snapshot the original values for the affected columns/rows before
# overwriting anything, so that reads are never affected by writes that
# already happened earlier in the loop.
# NOTE: use Value2 (not Value) - in this interop runtime, .Value getter
# does not reliably return the underlying scalar.
$original = @{}
foreach ($r in $sourceRows.Keys) {
    foreach ($c in $cols) {
        $addr = "$c$r"
        $original[$addr] = $ws.Range($addr).Value2
    }
}

foreach ($destRow in $sourceRows.Keys) {
    $srcRow = $sourceRows[$destRow]
    foreach ($c in $cols) {
        $srcAddr = "$c$srcRow"
        $destAddr = "$c$destRow"
        $ws.Range($destAddr).Value2 = $original[$srcAddr]
    }
}
